# "Drop in results from RMI script"
#
# 1. Remove the "Texas Notes" sheet (its content/notes are no longer
#    needed once real results are dropped in).
# 2. Update the Data sheet's C9:C11 "units sold" figures from placeholder
#    zeros to the real numbers from the RMI script, clearing the
#    yellow "needs data" highlight fill those cells carried while they
#    were still placeholders (reuse the neighbouring cells' formats so
#    no new style/fill entries are introduced).
# 3. The three BFoCSbQL-* sheets pull weighted averages from Data!C8:C11
#    etc., so they recalculate automatically once the Data values change.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. Drop the "Texas Notes" sheet -----------------------------------
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Delete()

# --- 2. Update Data!C9:C11 with the real figures -----------------------
$data = $wb.Worksheets.Item("Data")

# Borrow the (highlight-free) formatting already used by sibling cells
# in the same columns so we don't leave the old yellow "TBD" fill behind.
$data.Range("C8").Copy()
$data.Range("C9:C10").PasteSpecial(-4122)   # xlPasteFormats
$data.Range("D11").Copy()
$data.Range("C11").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = 0

$data.Range("C9").Value = 192000
$data.Range("C10").Value = 123000
$data.Range("C11").Value = 56000

# --- 3. Update cursor/selection bookmarks to match the refreshed file --
$data.Range("F30").Select()

$about = $wb.Worksheets.Item("About")
$about.Activate()
$about.Range("B17").Select()

$urban = $wb.Worksheets.Item("BFoCSbQL-urban-residential")
$urban.Range("A1").Select()

$rural = $wb.Worksheets.Item("BFoCSbQL-rural-residential")
$rural.Range("A1").Select()

$about.Activate()
